# Restore C10 (the "min" threshold for rule R30 on the Rules sheet) from
# 18 to 20, per revision #8105783ef8031ef75e825ea9b34a7e4c2903f2f0.TEST.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 20
